# Hopital Emergency Room PPT - slide 7 "Rectangle 2" shape:
# split "Dashboard 4: Consolidated View " into
# "Dashboard 4: " + "Keytakeaways" + " " (3 runs), keeping the
# existing run formatting (bold, accent5/lumMod 50%, green highlight,
# Arial Rounded MT Bold).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shp = $s.Shapes.Item(6)

$tr = $shp.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)

# Re-typing the paragraph text re-splits it into runs at the boundary
# where the original run ended, preserving that run's character
# formatting for the leading/trailing pieces and giving the newly
# typed middle word its own run.
$para1.Text = "Dashboard 4: Keytakeaways "
